# Remove the "Part 1 / Database Design and Implementation [65 Marks]" slide
# (presentation slide 5, p:sldId id="261") from the deck. All subsequent
# slides shift up by one; no other content changes.
$p = $ppt.ActivePresentation
$p.Slides.Item(5).Delete()
